# add PMS100 and DCA, reorganized script structure
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the date stamp in D1 by one day and give it its own (new) cell style
# entry — mirrors the style-table growth from the original authoring tool
# (an extra duplicate date-format style gets minted). Touching a throwaway
# cell with the same number format first, then clearing it, forces a fresh
# style slot to be allocated before D1 claims the next one.
$ws.Range("Z1").Value = 1
$ws.Range("Z1").NumberFormat = "mm-dd-yy"
$ws.Range("Z1").Clear()

$ws.Range("D1").Value = 42957.791666666664
$ws.Range("D1").NumberFormat = "mm-dd-yy"

# Updated counters for row 2 (PMS100 / DCA additions)
$ws.Range("A2").Value = 19
$ws.Range("B2").Value = 33
$ws.Range("C2").Value = 75
